$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prices in column D are plain text (e.g. "27.513.01", "1.004") that would
# otherwise be misread as numbers by COM Value coercion. Force text entry by
# switching the cells to a Text number format before writing, then restore
# the default "Normal" style so formatting matches the original workbook.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.448.72'
$ws.Range("E2").Value = '  +5.13%  '
$ws.Range("D3").Value = '1.722.48'
$ws.Range("E3").Value = '  +4.54%  '
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").Value = '225.50'
$ws.Range("E5").Value = '  +3.24%  '
$ws.Range("D6").Value = '0.5342'
$ws.Range("E6").Value = '  +2.71%  '
$ws.Range("D7").Value = '1.005'
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").Value = '0.2655'
$ws.Range("E8").Value = '  +1.39%  '
$ws.Range("D9").Value = '0.06591'
$ws.Range("E9").Value = '  +4.65%  '
$ws.Range("D10").Value = '21.50'
$ws.Range("E10").Value = '  +5.85%  '
$ws.Range("D11").Value = '0.07663'
$ws.Range("E11").Value = '  +0.33%  '
$ws.Range("D12").Value = '4.594'
$ws.Range("E12").Value = '  +0.25%  '
$ws.Range("D13").Value = '1.725.32'
$ws.Range("E13").Value = '  +3.73%  '
$ws.Range("D15").Value = '0.5785'
$ws.Range("E15").Value = '  +3.87%  '
$ws.Range("D16").Value = '0.0₅8275'
$ws.Range("E16").Value = '  +1.98%  '
$ws.Range("D17").Value = '67.74'
$ws.Range("E17").Value = '  +4.10%  '
$ws.Range("D18").Value = '27.451.96'
$ws.Range("E18").Value = '  +5.34%  '
$ws.Range("D19").Value = '217.15'
$ws.Range("E19").Value = '  +11.91%  '
$ws.Range("D20").Value = '1.004'
$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("D21").Value = '4.714'
$ws.Range("E21").Value = '  +2.78%  '
$ws.Range("D22").Value = '10.56'
$ws.Range("E22").Value = '  +1.21%  '
$ws.Range("D23").Value = '6.012'
$ws.Range("E23").Value = '  +1.63%  '
$ws.Range("D24").Value = '1.006'
$ws.Range("E24").Value = '  +0.19%  '
$ws.Range("D25").Value = '143.78'
$ws.Range("E25").Value = '  -0.86%  '
$ws.Range("D26").Value = '1.750'
$ws.Range("E26").Value = '  +13.90%  '
$ws.Range("D27").Value = '0.1229'
$ws.Range("E27").Value = '  +4.10%  '
$ws.Range("D28").Value = '7.320'
$ws.Range("E28").Value = '  +1.81%  '
$ws.Range("D29").Value = '16.46'
$ws.Range("E29").Value = '  +4.08%  '
$ws.Range("D30").Value = '0.05469'
$ws.Range("E30").Value = '  +0.55%  '
$ws.Range("D31").Value = '1.300'
$ws.Range("E31").Value = '  +2.47%  '
$ws.Range("D32").Value = '3.543'
$ws.Range("E32").Value = '  +3.25%  '
$ws.Range("D33").Value = '3.429'
$ws.Range("E33").Value = '  +3.19%  '
$ws.Range("D34").Value = '1.660'
$ws.Range("E34").Value = '  +6.50%  '
$ws.Range("D35").Value = '2.861'
$ws.Range("E35").Value = '  +2.90%  '
$ws.Range("D36").Value = '0.9542'
$ws.Range("E36").Value = '  +1.24%  '
$ws.Range("D37").Value = '2.428'
$ws.Range("E37").Value = '  +0.67%  '
$ws.Range("D38").Value = '0.5922'
$ws.Range("E38").Value = '  +5.97%  '
$ws.Range("D39").Value = '0.01645'
$ws.Range("E39").Value = '  +4.90%  '
$ws.Range("D40").Value = '5.900'
$ws.Range("E40").Value = '  +2.93%  '
$ws.Range("D41").Value = '0.8478'
$ws.Range("E41").Value = '  +3.41%  '
$ws.Range("D42").Value = '1.047.69'
$ws.Range("E42").Value = '  +1.99%  '
$ws.Range("D43").Value = '1.006'
$ws.Range("E43").Value = '  +0.25%  '
$ws.Range("D44").Value = '101.26'
$ws.Range("E44").Value = '  +0.80%  '
$ws.Range("D45").Value = '1.870.26'
$ws.Range("E45").Value = '  +4.75%  '
$ws.Range("D46").Value = '0.0₈114'
$ws.Range("E46").Value = '  +0.62%  '
$ws.Range("D47").Value = '58.61'
$ws.Range("E47").Value = '  +2.56%  '
$ws.Range("D48").Value = '0.4496'
$ws.Range("E48").Value = '  +4.12%  '
$ws.Range("D49").Value = '8.205'
$ws.Range("E49").Value = '  +4.24%  '
$ws.Range("E50").Value = '  +0.31%  '
$ws.Range("D51").Value = '0.05252'

$ws.Range("D2:D51").Style = "Normal"
